$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.09982999999999999
$ws.Range("H2").Value = 0.29949
$ws.Range("I2").Value = 0.06953924013247029
$ws.Range("J2").Value = 0.06953924013247029
$ws.Range("M2").Value = 8.252454666666667
$ws.Range("N2").Value = 24.757364
$ws.Range("O2").Value = 0.05349680956196952
$ws.Range("P2").Value = 0.05349680956196953
$ws.Range("Q2").Value = 0.8238425493733332
$ws.Range("R2").Value = 7.41458294436
$ws.Range("S2").Value = 0.003720127486450831
$ws.Range("T2").Value = 0.003720127486450832

$ws.Range("G3").Value = 0.09982999999999999
$ws.Range("H3").Value = 0.29949
$ws.Range("I3").Value = 0.06953924013247029
$ws.Range("J3").Value = 0.06953924013247029
$ws.Range("O3").Value = 0.5638948237978928
$ws.Range("P3").Value = 0.5638948237978929
$ws.Range("Q3").Value = 8.683892610043333
$ws.Range("R3").Value = 78.15503349039
$ws.Range("S3").Value = 0.03921281756153869
$ws.Range("T3").Value = 0.03921281756153869

$ws.Range("G4").Value = 0.09982999999999999
$ws.Range("H4").Value = 0.29949
$ws.Range("I4").Value = 0.06953924013247029
$ws.Range("J4").Value = 0.06953924013247029
$ws.Range("M4").Value = 57.81408433333333
$ws.Range("N4").Value = 173.442253
$ws.Range("O4").Value = 0.3747817085348802
$ws.Range("P4").Value = 0.3747817085348802
$ws.Range("Q4").Value = 5.771580038996666
$ws.Range("R4").Value = 51.94422035096999
$ws.Range("S4").Value = 0.02606203522706452
$ws.Range("T4").Value = 0.02606203522706453

$ws.Range("G5").Value = 0.09982999999999999
$ws.Range("H5").Value = 0.29949
$ws.Range("I5").Value = 0.06953924013247029
$ws.Range("J5").Value = 0.06953924013247029
$ws.Range("M5").Value = 1.207345666666667
$ws.Range("N5").Value = 3.622037
$ws.Range("O5").Value = 0.007826658105257385
$ws.Range("P5").Value = 0.007826658105257386
$ws.Range("Q5").Value = 0.1205293179033333
$ws.Range("R5").Value = 1.08476386113
$ws.Range("S5").Value = 0.0005442598574162382
$ws.Range("T5").Value = 0.0005442598574162383

$ws.Range("I6").Value = 0.4393303855760352
$ws.Range("J6").Value = 0.4393303855760352
$ws.Range("M6").Value = 8.252454666666667
$ws.Range("N6").Value = 24.757364
$ws.Range("O6").Value = 0.05349680956196952
$ws.Range("P6").Value = 0.05349680956196953
$ws.Range("Q6").Value = 5.204817656630222
$ws.Range("R6").Value = 46.84335890967201
$ws.Range("S6").Value = 0.0235027739719478
$ws.Range("T6").Value = 0.0235027739719478

$ws.Range("I7").Value = 0.4393303855760352
$ws.Range("J7").Value = 0.4393303855760352
$ws.Range("O7").Value = 0.5638948237978928
$ws.Range("P7").Value = 0.5638948237978929
$ws.Range("S7").Value = 0.2477361303634587
$ws.Range("T7").Value = 0.2477361303634587

$ws.Range("I8").Value = 0.4393303855760352
$ws.Range("J8").Value = 0.4393303855760352
$ws.Range("M8").Value = 57.81408433333333
$ws.Range("N8").Value = 173.442253
$ws.Range("O8").Value = 0.3747817085348802
$ws.Range("P8").Value = 0.3747817085348802
$ws.Range("Q8").Value = 36.46330444631045
$ws.Range("R8").Value = 328.169740016794
$ws.Range("S8").Value = 0.1646529925174741
$ws.Range("T8").Value = 0.1646529925174742

$ws.Range("I9").Value = 0.4393303855760352
$ws.Range("J9").Value = 0.4393303855760352
$ws.Range("M9").Value = 1.207345666666667
$ws.Range("N9").Value = 3.622037
$ws.Range("O9").Value = 0.007826658105257385
$ws.Range("P9").Value = 0.007826658105257386
$ws.Range("Q9").Value = 0.7614721070695556
$ws.Range("R9").Value = 6.853248963626001
$ws.Range("S9").Value = 0.003438488723154528
$ws.Range("T9").Value = 0.003438488723154529

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.705063
$ws.Range("H10").Value = 2.115189
$ws.Range("I10").Value = 0.4911303742914945
$ws.Range("J10").Value = 0.4911303742914945
$ws.Range("M10").Value = 8.252454666666667
$ws.Range("N10").Value = 24.757364
$ws.Range("O10").Value = 0.05349680956196952
$ws.Range("P10").Value = 0.05349680956196953
$ws.Range("Q10").Value = 5.818500444644
$ws.Range("R10").Value = 52.366504001796
$ws.Range("S10").Value = 0.02627390810357089
$ws.Range("T10").Value = 0.0262739081035709

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.705063
$ws.Range("H11").Value = 2.115189
$ws.Range("I11").Value = 0.4911303742914945
$ws.Range("J11").Value = 0.4911303742914945
$ws.Range("O11").Value = 0.5638948237978928
$ws.Range("P11").Value = 0.5638948237978929
$ws.Range("Q11").Value = 61.33117675363101
$ws.Range("R11").Value = 551.9805907826791
$ws.Range("S11").Value = 0.2769458758728954
$ws.Range("T11").Value = 0.2769458758728955

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.705063
$ws.Range("H12").Value = 2.115189
$ws.Range("I12").Value = 0.4911303742914945
$ws.Range("J12").Value = 0.4911303742914945
$ws.Range("M12").Value = 57.81408433333333
$ws.Range("N12").Value = 173.442253
$ws.Range("O12").Value = 0.3747817085348802
$ws.Range("P12").Value = 0.3747817085348802
$ws.Range("Q12").Value = 40.762571742313
$ws.Range("R12").Value = 366.863145680817
$ws.Range("S12").Value = 0.1840666807903415
$ws.Range("T12").Value = 0.1840666807903416

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.705063
$ws.Range("H13").Value = 2.115189
$ws.Range("I13").Value = 0.4911303742914945
$ws.Range("J13").Value = 0.4911303742914945
$ws.Range("M13").Value = 1.207345666666667
$ws.Range("N13").Value = 3.622037
$ws.Range("O13").Value = 0.007826658105257385
$ws.Range("P13").Value = 0.007826658105257386
$ws.Range("Q13").Value = 0.851254757777
$ws.Range("R13").Value = 7.661292819993
$ws.Range("S13").Value = 0.003843909524686619
$ws.Range("T13").Value = 0.00384390952468662
